# Add two new daily-attendance sheets ("2021-01-05" and "2021-01-06") at
# the end of the workbook, matching the layout used by the existing
# "2021-01-04" sheet (header row + data rows, header/column-A cells bold
# with a thin border, centered/top-aligned).

$wb = $excel.ActiveWorkbook

# Reference sheet whose formatting we reuse so no new style entries need
# to be invented (keeps the header + "Sr. No" column styling consistent
# with every other day-sheet in the workbook).
$styleSrc = $wb.Worksheets.Item("2021-01-04")

$headers = @("Sr. No","Name","Address","Job","Time-Stamp","SpO2_value","Heart-rate","Compensated","Ambient")

$rows0105 = @(
    @(1, "sachin",  "301/Sanskruti-1,Andheri, Mumbai",     "Software Engineer", "23:15:01", 97.94802944634111, 97.15533902413459,  "NA", "NA"),
    @(1, "sachin",  "301/Sanskruti-1,Andheri, Mumbai",     "Software Engineer", "23:17:25", 96.72572416450372, 67.29117181934262,  "NA", "NA"),
    @(1, "sachin",  "301/Sanskruti-1,Andheri, Mumbai",     "Software Engineer", "23:20:08", 97.01874773506313, 72.38625543168776,  "NA", "NA"),
    @(4, "dishant", "802/ Gunjan nagar/ Andheri , Mumbai", "Team Lead",         "23:22:31", 97.5138432365595,  62.76942121792798,  "NA", "NA"),
    @(1, "sachin",  "301/Sanskruti-1,Andheri, Mumbai",     "Software Engineer", "23:40:33", 97.62289986582616, 77.90709489973645,  "NA", "NA"),
    @(1, "sachin",  "301/Sanskruti-1,Andheri, Mumbai",     "Software Engineer", "23:45:52", 98.17296855471439, 89.27848258751145,  "NA", "NA"),
    @(1, "sachin",  "301/Sanskruti-1,Andheri, Mumbai",     "Software Engineer", "23:46:49", 97.34296827829047, 122.767066460571,   "NA", "NA"),
    @(1, "sachin",  "301/Sanskruti-1,Andheri, Mumbai",     "Software Engineer", "23:51:16", 97.77812564676016, 114.6530935674223,  "NA", "NA")
)

$rows0106 = @(
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "00:02:34", 97.67929655448826, 99.77017618367501, "NA", "NA"),
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "00:05:03", 97.24904921482502, 70.90286282809069, "NA", "NA"),
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "00:07:03", 97.96449387111554, 108.0425018446584, "NA", "NA"),
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "00:17:04", 97.03484786647245, 69.14380743344644, "NA", "NA")
)

# ---- Sheet "2021-01-05" ----
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws1 = $wb.Worksheets.Add([System.Type]::Missing, $after)
$ws1.Name = "2021-01-05"

for ($col = 1; $col -le $headers.Count; $col++) {
    $ws1.Cells.Item(1, $col).Value = $headers[$col - 1]
}
for ($r = 0; $r -lt $rows0105.Count; $r++) {
    $rowData = $rows0105[$r]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $ws1.Cells.Item($r + 2, $col).Value = $rowData[$col - 1]
    }
}

$styleSrc.Range("A1").Copy()
$ws1.Range("A1:I1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$ws1.Range("A2:A9").PasteSpecial(-4122)

# ---- Sheet "2021-01-06" ----
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $after)
$ws2.Name = "2021-01-06"

for ($col = 1; $col -le $headers.Count; $col++) {
    $ws2.Cells.Item(1, $col).Value = $headers[$col - 1]
}
for ($r = 0; $r -lt $rows0106.Count; $r++) {
    $rowData = $rows0106[$r]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $ws2.Cells.Item($r + 2, $col).Value = $rowData[$col - 1]
    }
}

$styleSrc.Range("A1").Copy()
$ws2.Range("A1:I1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$ws2.Range("A2:A5").PasteSpecial(-4122)

Write-Output "Added sheets 2021-01-05 and 2021-01-06"
